# feat: add 2022-Q1 data
#
# The workbook has 4 sheets: 2020-Q4, 2021-Q3, 2021-Q4, 总计 (a running
# summary sheet). This change:
#   1. Turns the old "总计" sheet (sheetId 4) into the new "2022-Q1"
#      per-fund holdings sheet.
#   2. Appends a brand-new "总计" sheet (sheetId 5) at the end with the
#      summary table, including a new row for 2022-Q1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Repurpose the old "总计" sheet as the "2022-Q1" holdings sheet.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(4)
$ws.Name = "2022-Q1"

# Paint the bold/bordered header style (already sitting on D1) onto the
# new header cells E1:H1, and the bold/bordered index-column style
# (already sitting on A2) onto the new index cells A3:A9, before the
# data is written.
$ws.Range("D1").Copy()
$ws.Range("E1:H1").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("A3:A9").PasteSpecial(-4122)

# Header row (plain text, never numeric-looking, so no special handling
# needed to keep it text).
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Columns B, D, E, F and G carry numeric-looking text (fund codes with
# leading zeros, and formatted decimal strings) that must stay text
# rather than being auto-coerced into numbers. Format as text first,
# write the value, then drop the formatting back to Normal so the
# cell ends up with the plain default style (matching the rest of the
# sheet) while keeping its stored type as text.
function Set-TextCell($addr, $text) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

function Set-Row($r, $code, $name, $size, $pos, $ratio, $value, $rank) {
    $ws.Range("A$r").Value = $r - 2
    Set-TextCell "B$r" $code
    $ws.Range("C$r").Value = $name
    Set-TextCell "D$r" $size
    Set-TextCell "E$r" $pos
    Set-TextCell "F$r" $ratio
    Set-TextCell "G$r" $value
    $ws.Range("H$r").Value = $rank
}

Set-Row 2 "011420" "广发全球科技三个月定期开放混合（QDII）人民币A" "34.16" "86.35" "4.78" "1.6328" 7
Set-Row 3 "011421" "广发全球科技三个月定期开放混合（QDII）美元A" "34.16" "86.35" "4.78" "1.6328" 7
Set-Row 4 "270023" "广发全球精选股票(QDII)" "25.53" "78.43" "5.63" "1.4373" 3
Set-Row 5 "000906" "广发全球精选股票(QDII)美元现汇" "25.53" "78.43" "5.63" "1.4373" 3
Set-Row 6 "011422" "广发全球科技三个月定期开放混合（QDII）人民币C" "6.53" "86.35" "4.78" "0.3121" 7
Set-Row 7 "011423" "广发全球科技三个月定期开放混合（QDII）美元C" "6.53" "86.35" "4.78" "0.3121" 7
Set-Row 8 "100055" "富国全球科技互联网股票(QDII)" "3.01" "70.87" "3.62" "0.1090" 7
Set-Row 9 "378006" "上投摩根全球新兴市场混合(QDII)" "0.46" "88.99" "6.84" "0.0315" 2

# ---------------------------------------------------------------------
# 2. Append a new "总计" summary sheet at the end of the workbook.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$totals = $wb.Worksheets.Add($null, $lastSheet)
$totals.Name = "总计"

# Apply the bold/bordered header style (row 1) and index-column style
# (column A) to match the rest of the workbook's sheets, sourcing the
# formatting from the sheet we just finished building.
$ws.Range("B1:D1").Copy()
$totals.Range("B1").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$totals.Range("A2:A5").PasteSpecial(-4122)

# Header row
$totals.Range("B1").Value = "日期"
$totals.Range("C1").Value = "持有数量(只)"
$totals.Range("D1").Value = "持有市值(亿元)"

function Set-TotalRow($r, $period, $count, $value) {
    $totals.Range("A$r").Value = $r - 2
    $totals.Range("B$r").Value = $period
    $totals.Range("C$r").Value = $count
    $totals.Range("D$r").Value = $value
}

Set-TotalRow 2 "2022-Q1" 8 6.9
Set-TotalRow 3 "2021-Q4" 4 3.77
Set-TotalRow 4 "2021-Q3" 3 0.12
Set-TotalRow 5 "2020-Q4" 9 1.4

# Restore the original active sheet/selection (sheet 1 was active before
# this edit and the diff doesn't touch that).
$firstSheet = $wb.Worksheets.Item(1)
[void]$firstSheet.Activate()
[void]$firstSheet.Range("A1").Select()
